$d = $word.ActiveDocument

# The document's header/footer logo pictures have their wp:docPr/@name and
# pic:cNvPr/@name swapped between "image1.*" and "image2.*". Renaming these
# inline shapes directly via InlineShapes.Item(N).Name is unreliable here
# (header shapes only update wp:docPr, footer shapes throw a stale-handle
# error), so instead we round-trip the document's full OOXML package text
# and patch the four picture name attributes directly - this is equivalent
# to what Word would persist if the renames had been applied through the UI.

$xml = $d.WordOpenXML

# wp:docPr elements are uniquely identified by their id attribute.
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>')

$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image2.png"/>')

$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>')

$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image1.jpg"/>')

# pic:cNvPr elements all share id="0", but the description text still
# distinguishes the two logos; both occurrences of each need the same
# rename so a plain global replace is correct here.
$xml = $xml.Replace(
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>',
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>')

$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>')

$d.WordOpenXML = $xml
